$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Units")
$ws.Range("AF2:AF6").Validation.Delete()
$ws.Range("AF8").Validation.Delete()
Write-Host "deleted"
